$d = $word.ActiveDocument

# Locate the paragraph that ends with "... para realizar o drag"
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*realizar o drag*") {
        $target = $p
    }
}

# Position right before the paragraph mark of that paragraph so the new
# text inherits that paragraph's run/paragraph formatting (pt-BR lang,
# no extra color) instead of the formatting of whatever follows it.
$insertPos = $target.Range.End - 1
$rng = $d.Range($insertPos, $insertPos)

$rng.InsertAfter("`r*Erro retornado quando orderList é null no reducer`r*BackgroundColor aparecendo em qualquer click, fazer com que apareça somente quando for acontecer o drag ou simplesmente fazer com que suma caso não ocorra o drag")
